$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: reported_product_problems/array -> number_tobacco_products/number
$ws.Range("B4").Value = "number_tobacco_products"
$ws.Range("C4").Value = "number"
$ws.Range("D4").Value = "System-calculated number of Tobacco Product Problems reported, displayed as a whole number, ≥ 0."
$ws.Rows(4).RowHeight = 34

# Row 5: reported_health_problems/array -> tobacco_products/array
$ws.Range("B5").Value = "tobacco_products"
$ws.Range("C5").Value = "array"
$ws.Range("D5").Value = "Text reflecting the SRP tobacco Product Type selected by the reporter."
$ws.Rows(5).RowHeight = 17

# Row 6: product_type/array -> number_health_problems/number
$ws.Range("B6").Value = "number_health_problems"
$ws.Range("C6").Value = "number"
$ws.Range("D6").Value = "System-calculated number of Health Problems (i.e., MedDRA terms selected from a standardized list of symptoms, signs, diagnoses and outcomes) reported, displayed as a whole number, ≥0."
$ws.Rows(6).RowHeight = 51

# Row 7: number_tobacco_products/number -> reported_health_problems/array
$ws.Range("B7").Value = "reported_health_problems"
$ws.Range("C7").Value = "array"
$ws.Range("D7").Value = "Text reflecting the MedDRA terms selected by the reporter."
$ws.Rows(7).RowHeight = 17

# Row 8: number_product_problems/number -> nonuser_affected/string
$ws.Range("B8").Value = "nonuser_affected"
$ws.Range("C8").Value = "string"
$ws.Range("D8").Value = "Displays text reflecting the response to this optional question (2017 - 12/14/2018) or required question (12/15/2018 onward) as “No information provided” if not answered, or Yes/No."
$ws.Rows(8).RowHeight = 51

# Row 9: number_health_problems/number -> number_product_problems/number
$ws.Range("B9").Value = "number_product_problems"
$ws.Range("C9").Value = "number"
$ws.Range("D9").Value = "System-calculated number of categorical Product Problems reported, displayed as a whole number, ≥0."
$ws.Rows(9).RowHeight = 34

# Row 10: nonuser_affected/string -> reported_product_problems/array
$ws.Range("B10").Value = "reported_product_problems"
$ws.Range("C10").Value = "array"
$ws.Range("D10").Value = "Text reflecting the SRP categorical list of values."
$ws.Rows(10).RowHeight = 17

# Update the active selection to match the new view state
$ws.Range("B12").Select()
